# #17 checking menu mockup
# - Re-stamp the cached "datetimeFigureOut" field text (11/29/2021 -> 11/28/2021)
#   on the slide master and every slide layout's Date Placeholder.
# - Fix up the "checking menu" mockup textboxes on slide 3 (typos /
#   wording tweaks) and resize a handful of the red mock-up labels to
#   match their new (slightly longer/shorter) text.

$p = $ppt.ActivePresentation

$oldDate = "11/29/2021"
$newDate = "11/28/2021"

# --- Slide master: Date Placeholder -----------------------------------
for ($i = 1; $i -le $p.SlideMaster.Shapes.Count; $i++) {
    $sh = $p.SlideMaster.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every slide layout: Date Placeholder ------------------------------
for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $lyt = $p.SlideMaster.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lyt.Shapes.Count; $i++) {
        $sh = $lyt.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- Slide 3: mockup labels ---------------------------------------------
$slide = $p.Slides.Item(3)

# "<Even>" -> "<Event>", widen the textbox to fit the new text
$evenBox = $slide.Shapes.Item(13)
$evenBox.TextFrame.TextRange.Text = "<Event>"
$evenBox.Width = 74.37433270866141

# "<FindEvent>" -> "<allevent>", narrow the textbox to fit the new text
$findEventBox = $slide.Shapes.Item(14)
$findEventBox.TextFrame.TextRange.Text = "<allevent>"
$findEventBox.Width = 93.2367746535433

# "<category>" -> "<categories>", widen the textbox to fit the new text
$categoryBox = $slide.Shapes.Item(16)
$categoryBox.TextFrame.TextRange.Text = "<categories>"
$categoryBox.Width = 110.9227569055118

# "// category gets userdata to app  " -> "// categories gets userdata to app  "
$cmt1 = $slide.Shapes.Item(26)
$cmt1.TextFrame.TextRange.Characters(4, 9).Text = "categories "
$cmt1.Width = 211.03095288188976
$cmt1.Height = 23.896142732283465

# "// findevent gets userdata to app  " -> "// allevent gets userdata to app  "
$cmt2 = $slide.Shapes.Item(27)
$cmt2.TextFrame.TextRange.Characters(4, 10).Text = "allevent "
$cmt2.Width = 198.06055518110236
$cmt2.Height = 23.896142732283465

# "// event gets userdata to app  " -> "// envent gets userdata to app  "
$cmt3 = $slide.Shapes.Item(28)
$cmt3.TextFrame.TextRange.Characters(4, 6).Text = "envent "
$cmt3.Width = 192.21393700787402
$cmt3.Height = 23.896142732283465
